$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.142.76"
$ws.Range("E2").Value = "  -0.32%  "

$ws.Range("D3").Value = "'1.668.72"
$ws.Range("E3").Value = "  -0.86%  "

$ws.Range("E4").Value = "  -0.52%  "

$ws.Range("D5").Value = "'209.45"
$ws.Range("E5").Value = "  -3.43%  "

$ws.Range("D6").Value = "'0.5230"
$ws.Range("E6").Value = "  -1.16%  "

$ws.Range("E7").Value = "  -0.52%  "

$ws.Range("D8").Value = "'0.2624"
$ws.Range("E8").Value = "  -2.91%  "

$ws.Range("D9").Value = "'0.06341"
$ws.Range("E9").Value = "  -0.53%  "

$ws.Range("D10").Value = "'21.19"
$ws.Range("E10").Value = "  -1.76%  "

$ws.Range("D11").Value = "'0.07527"
$ws.Range("E11").Value = "  -1.83%  "

$ws.Range("D12").Value = "'1.679.44"
$ws.Range("E12").Value = "  -1.25%  "

$ws.Range("D13").Value = "'4.445"
$ws.Range("E13").Value = "  -1.58%  "

$ws.Range("D14").Value = "'0.5499"
$ws.Range("E14").Value = "  -4.65%  "

$ws.Range("D15").Value = "'66.41"
$ws.Range("E15").Value = "  +0.21%  "

$ws.Range("D16").Value = "'0.000007964"
$ws.Range("E16").Value = "  -4.44%  "

$ws.Range("D17").Value = "'26.155.98"
$ws.Range("E17").Value = "  -0.44%  "

$ws.Range("D18").Value = "'1.003"
$ws.Range("E18").Value = "  -0.52%  "

$ws.Range("D19").Value = "'4.749"
$ws.Range("E19").Value = "  -2.53%  "

$ws.Range("D20").Value = "'186.61"
$ws.Range("E20").Value = "  -1.59%  "

$ws.Range("D21").Value = "'10.30"
$ws.Range("E21").Value = "  -4.84%  "

$ws.Range("D22").Value = "'6.185"
$ws.Range("E22").Value = "  -0.70%  "

$ws.Range("E23").Value = "  -0.51%  "

$ws.Range("D24").Value = "'149.46"
$ws.Range("E24").Value = "  +0.46%  "

$ws.Range("D25").Value = "'0.1250"
$ws.Range("E25").Value = "  -1.17%  "

$ws.Range("D26").Value = "'7.506"
$ws.Range("E26").Value = "  -3.74%  "

$ws.Range("D27").Value = "'15.86"
$ws.Range("E27").Value = "  +0.96%  "

$ws.Range("D28").Value = "'0.06425"
$ws.Range("E28").Value = "  +3.01%  "

$ws.Range("D29").Value = "'1.354"
$ws.Range("E29").Value = "  -1.42%  "

$ws.Range("D30").Value = "'1.275"
$ws.Range("E30").Value = "  -3.40%  "

$ws.Range("D31").Value = "'3.510"
$ws.Range("E31").Value = "  -1.78%  "

$ws.Range("D32").Value = "'3.415"
$ws.Range("E32").Value = "  -4.09%  "

$ws.Range("D33").Value = "'1.646"
$ws.Range("E33").Value = "  -2.37%  "

$ws.Range("D34").Value = "'1.006"
$ws.Range("E34").Value = "  -1.71%  "

$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").Value = "'0.6024"
$ws.Range("E35").Value = "  -1.99%  "

$ws.Range("B36").Value = "HuobiToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D36").Value = "'2.407"
$ws.Range("E36").Value = "  -0.79%  "

$ws.Range("D37").Value = "'2.743"
$ws.Range("E37").Value = "  -0.66%  "

$ws.Range("D38").Value = "'1.112.94"
$ws.Range("E38").Value = "  +0.73%  "

$ws.Range("D39").Value = "'6.139"
$ws.Range("E39").Value = "  -1.02%  "

$ws.Range("D40").Value = "'0.01616"
$ws.Range("E40").Value = "  -0.34%  "

$ws.Range("D41").Value = "'0.8679"
$ws.Range("E41").Value = "  -2.97%  "

$ws.Range("E42").Value = "  -0.81%  "

$ws.Range("D43").Value = "'100.31"
$ws.Range("E43").Value = "  -0.21%  "

$ws.Range("D44").Value = "'1.821.82"
$ws.Range("E44").Value = "  -0.71%  "

$ws.Range("E45").Value = "  -1.05%  "

$ws.Range("D46").Value = "'55.46"
$ws.Range("E46").Value = "  -3.55%  "

$ws.Range("D47").Value = "'1.003"
$ws.Range("E47").Value = "  +0.11%  "

$ws.Range("D48").Value = "'8.041"
$ws.Range("E48").Value = "  -0.67%  "

$ws.Range("D49").Value = "'0.05231"
$ws.Range("E49").Value = "  -0.82%  "

$ws.Range("D50").Value = "'0.4245"
$ws.Range("E50").Value = "  -1.04%  "

$ws.Range("D51").Value = "'5.929"
$ws.Range("E51").Value = "  -1.52%  "
